$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1146.3636
$ws.Range("I112").Value = 650
$ws.Range("J112").Value = 1256.6666
$ws.Range("K112").Value = 1950
$ws.Range("L112").Value = 3769.9998
$ws.Range("M112").Value = -842
$ws.Range("N112").Value = -5985.9998
$ws.Range("H116").Value = 2900.375
$ws.Range("I116").Value = 2565.6667
$ws.Range("J116").Value = 3101.2
$ws.Range("K116").Value = 2565.6667
$ws.Range("L116").Value = 3101.2
$ws.Range("N116").Value = -9985.200000000001
$ws.Range("M116").Value = 876.3332999999998
$ws.Range("H118").Value = 710.25
$ws.Range("I118").Value = 360.42856
$ws.Range("J118").Value = 1200
$ws.Range("K118").Value = 1081.28568
$ws.Range("L118").Value = 3600
$ws.Range("M118").Value = 575.71432
$ws.Range("N118").Value = -6914
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1726.25
$ws.Range("I61").Value = 1707.0667
$ws.Range("J61").Value = 2014
$ws.Range("K61").Value = 1707.0667
$ws.Range("L61").Value = 2014
$ws.Range("M61").Value = -1495.0667
$ws.Range("N61").Value = -2438
$ws.Range("H63").Value = 3544.6155
$ws.Range("I63").Value = 2180
$ws.Range("K63").Value = 2180
$ws.Range("M63").Value = -1494
$ws.Range("H66").Value = 3544.6155
$ws.Range("I66").Value = 2180
$ws.Range("K66").Value = 10900
$ws.Range("M66").Value = -7468
$ws.Range("H82").Value = 27500
$ws.Range("J82").Value = 27500
$ws.Range("L82").Value = 27500
$ws.Range("N82").Value = -28222
$ws.Range("H85").Value = 27500
$ws.Range("J85").Value = 27500
$ws.Range("L85").Value = 27500
$ws.Range("N85").Value = -29996
$ws.Range("H125").Value = 23311.818
$ws.Range("J125").Value = 23311.818
$ws.Range("L125").Value = 23311.818
$ws.Range("N125").Value = -33151.818
$ws.Range("H136").Value = 1726.25
$ws.Range("I136").Value = 1707.0667
$ws.Range("J136").Value = 2014
$ws.Range("K136").Value = 5121.2001
$ws.Range("L136").Value = 6042
$ws.Range("M136").Value = -2571.2001
$ws.Range("N136").Value = -11142
$ws.Range("H139").Value = 172333.33
$ws.Range("J139").Value = 172333.33
$ws.Range("L139").Value = 172333.33
$ws.Range("N139").Value = -182613.33
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1725.9736
$ws.Range("I86").Value = 1652.3529
$ws.Range("J86").Value = 1785.5714
$ws.Range("K86").Value = 1652.3529
$ws.Range("L86").Value = 1785.5714
$ws.Range("M86").Value = -529.3529000000001
$ws.Range("N86").Value = -4031.5714
$ws.Range("H89").Value = 1725.9736
$ws.Range("I89").Value = 1652.3529
$ws.Range("J89").Value = 1785.5714
$ws.Range("K89").Value = 8261.764500000001
$ws.Range("L89").Value = 8927.857
$ws.Range("M89").Value = -2645.764500000001
$ws.Range("N89").Value = -20159.857
$ws.Range("H103").Value = 35000.332
$ws.Range("J103").Value = 35000.332
$ws.Range("L103").Value = 35000.332
$ws.Range("N103").Value = -37344.332
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents() | Out-Null
$ws.Range("H62").Value = 2406.5625
$ws.Range("I62").Value = 2304.4546
$ws.Range("J62").Value = 2631.2
$ws.Range("K62").Value = 2304.4546
$ws.Range("L62").Value = 2631.2
$ws.Range("M62").Value = -1680.4546
$ws.Range("N62").Value = -3879.2
$ws.Range("H65").Value = 2406.5625
$ws.Range("I65").Value = 2304.4546
$ws.Range("J65").Value = 2631.2
$ws.Range("K65").Value = 11522.273
$ws.Range("L65").Value = 13156
$ws.Range("M65").Value = -8402.273000000001
$ws.Range("N65").Value = -19396
$ws.Range("H68").Value = 29500
$ws.Range("J68").Value = 29500
$ws.Range("L68").Value = 29500
$ws.Range("N68").Value = -30998
$ws.Range("H71").Value = 29500
$ws.Range("J71").Value = 29500
$ws.Range("L71").Value = 88500
$ws.Range("N71").Value = -95988
$ws.Range("H74").Value = 21250
$ws.Range("I74").Value = 15000
$ws.Range("J74").Value = 23333.334
$ws.Range("K74").Value = 15000
$ws.Range("L74").Value = 23333.334
$ws.Range("N74").Value = -25081.334
$ws.Range("M74").Value = -14126
$ws.Range("H77").Value = 21250
$ws.Range("I77").Value = 15000
$ws.Range("J77").Value = 23333.334
$ws.Range("K77").Value = 45000
$ws.Range("L77").Value = 70000.00199999999
$ws.Range("N77").Value = -78736.00199999999
$ws.Range("M77").Value = -40632
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1021.5294
$ws.Range("I5").Value = 415.5
$ws.Range("K5").Value = 1246.5
$ws.Range("M5").Value = -1134.5
$ws.Range("H23").Value = 157
$ws.Range("I23").Value = 56
$ws.Range("J23").Value = 229.14285
$ws.Range("K23").Value = 168
$ws.Range("L23").Value = 687.4285500000001
$ws.Range("M23").Value = 67
$ws.Range("N23").Value = -1157.42855
$ws.Range("H26").Value = 648.6667
$ws.Range("J26").Value = 913
$ws.Range("L26").Value = 2739
$ws.Range("N26").Value = -3315
$ws.Range("H113").Value = 1378034.8
$ws.Range("I113").Value = 2331615.2
$ws.Range("J113").Value = 640.7778
$ws.Range("K113").Value = 6994845.600000001
$ws.Range("L113").Value = 1922.3334
$ws.Range("M113").Value = -6992675.600000001
$ws.Range("N113").Value = -6262.3334
$ws.Range("H116").Value = 975
$ws.Range("I116").Value = 975
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2925
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 517
$ws.Range("N116").ClearContents() | Out-Null
$ws.Range("H124").Value = 3386.5
$ws.Range("J124").Value = 3492
$ws.Range("L124").Value = 10476
$ws.Range("N124").Value = -20296
$ws.Range("H135").Value = 1021.5294
$ws.Range("I135").Value = 415.5
$ws.Range("K135").Value = 3739.5
$ws.Range("M135").Value = -1204.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1585.5333
$ws.Range("I7").Value = 1521.7693
$ws.Range("K7").Value = 1521.7693
$ws.Range("M7").Value = -1409.7693
$ws.Range("H22").Value = 42465.875
$ws.Range("I22").Value = 250420.25
$ws.Range("J22").Value = 875
$ws.Range("K22").Value = 250420.25
$ws.Range("L22").Value = 875
$ws.Range("M22").Value = -250125.25
$ws.Range("N22").Value = -1465
$ws.Range("H27").Value = 42465.875
$ws.Range("I27").Value = 250420.25
$ws.Range("J27").Value = 875
$ws.Range("K27").Value = 250420.25
$ws.Range("L27").Value = 875
$ws.Range("M27").Value = -250313.25
$ws.Range("N27").Value = -1089
$ws.Range("H46").Value = 932.8461
$ws.Range("J46").Value = 783.06665
$ws.Range("L46").Value = 783.06665
$ws.Range("N46").Value = -1159.06665
$ws.Range("H126").Value = 1585.5333
$ws.Range("I126").Value = 1521.7693
$ws.Range("K126").Value = 4565.3079
$ws.Range("M126").Value = -2095.3079
$ws.Range("H137").Value = 66888.89
$ws.Range("I137").Value = 150000
$ws.Range("J137").Value = 56500
$ws.Range("K137").Value = 150000
$ws.Range("L137").Value = 56500
$ws.Range("M137").Value = -144900
$ws.Range("N137").Value = -66700
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 29286.666
$ws.Range("I82").Value = 29860
$ws.Range("J82").Value = 29000
$ws.Range("K82").Value = 29860
$ws.Range("L82").Value = 29000
$ws.Range("N82").Value = -29766
$ws.Range("M82").Value = -29477
$ws.Range("H85").Value = 29286.666
$ws.Range("I85").Value = 29860
$ws.Range("J85").Value = 29000
$ws.Range("K85").Value = 29860
$ws.Range("L85").Value = 29000
$ws.Range("N85").Value = -31652
$ws.Range("M85").Value = -28534

Write-Output "Applied 201 cell updates across 7 sheets"